# ---------------------------------------------------------------------------
# Weekly CompStat refresh: bump the report Volume/date-range header text and
# reload the precinct crime-complaint table (rows 15-29) with the new weekly
# figures. Mirrors the source system's "new crime data collected" refresh.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Volume 29   Number  42" -> "...  43" ------------------------
$a8 = $ws.Range("A8")
$a8.Characters(21, 2).Text = "43"
$ws.Range("A8").Characters(21, 2).Font.Name = "Andale WT"
$ws.Range("A8").Characters(21, 2).Font.Size = 10

# --- Header: week-covering date range shifts forward by one week -----------
$c9 = $ws.Range("C9")
$c9.Characters(27, 10).Text = "10/24/2022"
$ws.Range("C9").Characters(27, 10).Font.Name = "Andale WT"
$ws.Range("C9").Characters(27, 10).Font.Size = 10
$ws.Range("C9").Characters(48, 10).Text = "10/30/2022"
$ws.Range("C9").Characters(48, 10).Font.Name = "Andale WT"
$ws.Range("C9").Characters(48, 10).Font.Size = 10

# --- Crime-complaint grid: refreshed weekly / 28-day / YTD / 2-yr figures --
# Row 15
$ws.Range("F15").Value = 3
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 12
$ws.Range("J15").Value = 37
$ws.Range("K15").Value = -67.567567567567
$ws.Range("L15").Value = -62.5
$ws.Range("M15").Value = -25
$ws.Range("N15").Value = -63.636363636363
# Row 16
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 1
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = 116.666666666667
$ws.Range("I16").Value = 155
$ws.Range("J16").Value = 114
$ws.Range("K16").Value = 35.964912280701
$ws.Range("L16").Value = 6.896551724137
$ws.Range("M16").Value = -41.06463878327
$ws.Range("N16").Value = -81.997677119628
# Row 17
$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 80
$ws.Range("F17").Value = 26
$ws.Range("H17").Value = -23.529411764705
$ws.Range("I17").Value = 350
$ws.Range("J17").Value = 346
$ws.Range("K17").Value = 1.156069364161
$ws.Range("L17").Value = -3.314917127071
$ws.Range("M17").Value = 48.305084745762
$ws.Range("N17").Value = -0.2849002849
# Row 18
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 18
$ws.Range("H18").Value = 100
$ws.Range("I18").Value = 173
$ws.Range("J18").Value = 119
$ws.Range("K18").Value = 45.378151260504
$ws.Range("L18").Value = 2.366863905325
$ws.Range("M18").Value = -45.425867507886
$ws.Range("N18").Value = -87.118391660461
# Row 19
$ws.Range("C19").Value = 12
$ws.Range("E19").Value = -7.692307692307
$ws.Range("F19").Value = 38
$ws.Range("G19").Value = 46
$ws.Range("H19").Value = -17.391304347826
$ws.Range("I19").Value = 523
$ws.Range("J19").Value = 379
$ws.Range("K19").Value = 37.994722955145
$ws.Range("L19").Value = 3.155818540433
$ws.Range("M19").Value = 39.466666666666
$ws.Range("N19").Value = 4.39121756487
# Row 20
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 200
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = 191.666666666667
$ws.Range("I20").Value = 227
$ws.Range("J20").Value = 119
$ws.Range("K20").Value = 90.756302521008
$ws.Range("L20").Value = -12.015503875969
$ws.Range("M20").Value = -29.938271604938
$ws.Range("N20").Value = -91.973125884017
# Row 21
$ws.Range("C21").Value = 35
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = 34.615384615384
$ws.Range("F21").Value = 134
$ws.Range("G21").Value = 110
$ws.Range("H21").Value = 21.818181818181
$ws.Range("I21").Value = 1445
$ws.Range("J21").Value = 1118
$ws.Range("K21").Value = 29.248658318425
$ws.Range("L21").Value = -2.232746955345
$ws.Range("M21").Value = -6.411917098445
$ws.Range("N21").Value = -75.661108303857
# Row 24
$ws.Range("C24").Value = 39
$ws.Range("D24").Value = 34
$ws.Range("E24").Value = 14.705882352941
$ws.Range("F24").Value = 117
$ws.Range("G24").Value = 100
$ws.Range("H24").Value = 17
$ws.Range("I24").Value = 1273
$ws.Range("J24").Value = 826
$ws.Range("K24").Value = 54.11622276029
$ws.Range("L24").Value = 50.295159386068
$ws.Range("M24").Value = 80.311614730878
# Row 25
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = 33.333333333333
$ws.Range("F25").Value = 55
$ws.Range("G25").Value = 51
$ws.Range("H25").Value = 7.843137254901
$ws.Range("I25").Value = 562
$ws.Range("J25").Value = 422
$ws.Range("K25").Value = 33.175355450237
$ws.Range("L25").Value = 26.008968609865
$ws.Range("M25").Value = 9.551656920077
# Row 26
$ws.Range("C26").Value = 1
$ws.Range("E26").Value = 0
$ws.Range("I26").Value = 27
$ws.Range("J26").Value = 46
$ws.Range("K26").Value = -41.304347826087
$ws.Range("L26").Value = -30.76923076923
# Row 27
$ws.Range("F27").Value = 8
$ws.Range("H27").Value = 166.666666666667
$ws.Range("I27").Value = 54
$ws.Range("K27").Value = 28.571428571428
$ws.Range("L27").Value = 22.727272727272
# Row 28
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 0
$ws.Range("M28").Value = -47.619047619047
$ws.Range("N28").Value = -67.647058823529
# Row 29
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = 0
$ws.Range("M29").Value = -33.333333333333
$ws.Range("N29").Value = -70

# --- Cells that flip between a numeric figure and the "N/A" (0 / ***.*) ----
# --- placeholder pair used when a category has no prior-year baseline ------

# Row 18: D18/E18 go from "N/A" placeholders to real figures
$ws.Range("D18").NumberFormat = "#,##0"
$ws.Range("D18").Value = 3
$ws.Range("E18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E18").Value = -33.333333333333

# Row 27: C27 gains a real figure; D27/E27 become the new "N/A" placeholders
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("C27").Value = 6

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("D27").ClearFormats()
$ws.Range("D27").Font.Name = "Andale WT"
$ws.Range("D27").Font.Size = 10
$ws.Range("D27").HorizontalAlignment = -4152
$ws.Range("D27").VerticalAlignment = -4108

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("E27").ClearFormats()
$ws.Range("E27").Font.Name = "Andale WT"
$ws.Range("E27").Font.Size = 10
$ws.Range("E27").HorizontalAlignment = -4152
$ws.Range("E27").VerticalAlignment = -4108
